$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5971993803977966
$ws.Range("B1").Value = 1.490510106086731
$ws.Range("C1").Value = 5.93446159362793
$ws.Range("D1").Value = 2.040632009506226
$ws.Range("E1").Value = 1.51197350025177
